$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(29, 0.28, 4808.447313954, 34.602017566, 4773.845296388, 28.35622475466667, 1687.40019971, 24.411650662, 1662.98854905, 29.71434492466667),
    @(30, 0.29, 4855.412432398, 35.855350794, 4819.557081604, 91.05551909399999, 1691.14584605, 26.621560398, 1664.524285652, 95.16820697600001)
)

foreach ($rowData in $data) {
    $r = $rowData[0]
    for ($col = 1; $col -le 9; $col++) {
        $ws.Cells.Item($r, $col).Value = $rowData[$col]
    }
}
